$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.815.27'
$ws.Range('E2').Value = '  +0.41%  '

# Row 3
$ws.Range('D3').Value = '3.329.15'
$ws.Range('E3').Value = '  +2.39%  '

# Row 4
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '605.03'
$ws.Range('E5').Value = '  +1.07%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.95'
$ws.Range('E6').Value = '  +0.89%  '

# Row 7
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').Value = '3.327.83'
$ws.Range('E8').Value = '  +2.14%  '

# Row 9
$ws.Range('E9').Value = '  -0.48%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  +1.14%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.54'
$ws.Range('E11').Value = '  +2.28%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.471'
$ws.Range('E12').Value = '  +0.64%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  -0.84%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '35.07'
$ws.Range('E14').Value = '  +0.67%  '

# Row 15
$ws.Range('D15').Value = '3.873.44'
$ws.Range('E15').Value = '  +2.53%  '

# Row 16
$ws.Range('E16').Value = '  +0.06%  '

# Row 17
$ws.Range('D17').Value = '3.323.26'
$ws.Range('E17').Value = '  +2.97%  '

# Row 18
$ws.Range('D18').Value = '63.902.76'
$ws.Range('E18').Value = '  +0.69%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.88'
$ws.Range('E19').Value = '  +0.80%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '480.82'
$ws.Range('E20').Value = '  -0.13%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.12'
$ws.Range('E21').Value = '  -0.96%  '

# Row 22
$ws.Range('E22').Value = '  +1.56%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.98'
$ws.Range('E23').Value = '  +0.22%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.78'
$ws.Range('E24').Value = '  +3.07%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.93'
$ws.Range('E25').Value = '  +0.55%  '

# Row 26
$ws.Range('E26').Value = '  -0.02%  '

# Row 27
$ws.Range('E27').Value = '  +1.15%  '

# Row 28
$ws.Range('E28').Value = '  +1.57%  '

# Row 29
$ws.Range('E29').Value = '  +0.01%  '

# Row 30
$ws.Range('E30').Value = '  -2.32%  '

# Row 31
$ws.Range('E31').Value = '  +0.97%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '28.96'
$ws.Range('E32').Value = '  +4.22%  '

# Row 33
$ws.Range('E33').Value = '  -1.88%  '

# Row 34
$ws.Range('E34').Value = '  -1.98%  '

# Row 35
$ws.Range('E35').Value = '  +0.68%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.10'
$ws.Range('E36').Value = '  +2.79%  '

# Row 37
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0749'
$ws.Range('E37').Value = '  +2.28%  '

# Row 38
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '52.39'
$ws.Range('E38').Value = '  -1.19%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0402'
$ws.Range('E39').Value = '  +1.78%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '434.53'
$ws.Range('E40').Value = '  +1.33%  '

# Row 41
$ws.Range('D41').Value = '3.098.25'
$ws.Range('E41').Value = '  +3.37%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.78'
$ws.Range('E42').Value = '  -1.17%  '

# Row 43
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.118'
$ws.Range('E43').Value = '  +5.16%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.36'
$ws.Range('E44').Value = '  -0.98%  '

# Row 45
$ws.Range('E45').Value = '  -0.71%  '

# Row 46
$ws.Range('E46').Value = '  +3.94%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '36.99'
$ws.Range('E47').Value = '  +12.06%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '26.50'
$ws.Range('E48').Value = '  +1.36%  '

# Row 49
$ws.Range('E49').Value = '  -0.02%  '

# Row 50
$ws.Range('E50').Value = '  -0.11%  '

# Row 51
$ws.Range('E51').Value = '  -0.64%  '
